# Summarizing what we did with the simulations
#
# The last list item ("Negative Skewed" -> a=1, b=4) currently holds its
# text split across two runs ("a" and "=1, b=4") with the _GoBack bookmark
# sitting in between them. The edit:
#   1. Merges that text into a single run reading "a=1, b=4".
#   2. Moves the _GoBack bookmark into its own (now empty) list paragraph
#      at the top list level (ilvl 0).
#   3. Adds a new top-level (ilvl 0) list paragraph with a new sentence
#      about the WinGen/theta comparison.

$d = $word.ActiveDocument

# Paragraph 11 is "a" + bookmark + "=1, b=4"  ->  text should read "a=1, b=4"
$p = $d.Paragraphs.Item(11)

# Replace the first run's text ("a") with the fully merged text.
$rng = $p.Range
$rng.Find.Execute("a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "a=1, b=4"

# Re-fetch the paragraph and drop the now-duplicated trailing "=1, b=4"
# that used to live in the second run (right after the merged text).
$p = $d.Paragraphs.Item(11)
$tailStart = $p.Range.Start + 8
$tailRng = $d.Range($tailStart, $p.Range.End)
$tailRng.Delete()

# Remove the old bookmark; it will be re-created in its own paragraph below.
$d.Bookmarks("_GoBack").Delete()

# Build the two new list paragraphs (bookmark-only paragraph, then the new
# sentence) and splice them in right after paragraph 11, both at ilvl 0.
$p = $d.Paragraphs.Item(11)
$insertionPoint = $d.Range($p.Range.End, $p.Range.End)

$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">We found a relationship between the slope and the ability (thetas). In </w:t></w:r>
<w:proofErr w:type="spellStart"/><w:r><w:t>WinGen</w:t></w:r><w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> the thetas are generated differently, with more variability.</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>The more variability, the smaller the slope.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newParagraphsXml) | Out-Null
